$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted before the existing row 205 in the
# Pepino ensalada dataset. Inserting the row shifts every subsequent record
# (old rows 205..282) down by one (new rows 206..283), exactly matching the
# commit diff. Excel's normal "Insert" semantics (shift down, carry the
# formatting of the row above) are used, matching how row 205's D column
# keeps style index 2 after the insert.
$ws.Rows.Item(205).Insert()

# Populate the freshly inserted row 205 with the new record's data.
$ws.Range("A205").Value = 7
$ws.Range("B205").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C205").Value = "Ñuble"
$ws.Range("D205").Value = 44924
$ws.Range("E205").Value = 16
$ws.Range("F205").Value = 100112043
$ws.Range("G205").Value = "Pepino ensalada"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 120
$ws.Range("K205").Value = 13000
$ws.Range("L205").Value = 14000
$ws.Range("M205").Value = 13500
$ws.Range("N205").Value = "$/caja 80 unidades"
$ws.Range("O205").Value = "Región del Maule"
$ws.Range("P205").Value = 169
$ws.Range("Q205").Value = 80
$ws.Range("R205").Value = "Hortaliza"
